$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A4 with the new combined tuple-style strings
$ws.Range("A2").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A3").Value = "('Sliver', ['Token Creature — Sliver', '1/1'])"
$ws.Range("A4").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"

# Remove the now-unused rows 5 through 14 (shifting cells up)
$ws.Range("A5:A14").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
